# Applies the weekly update to the "Fruta, Macroferia Regional de Talca - Chirimoya" sheet:
#  - inserts one new record (2021-08-13 / 44421) above the current row 5, pushing the
#    existing rows 5-18 down to rows 6-19
#  - appends two new records (both dated 2021-09-22 / 44461) as new rows 20 and 21

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common (constant) values shared by every data row in this table.
$mercadoId = 5
$mercado   = "Macroferia Regional de Talca"
$region    = "Maule"
$codreg    = 7
$tipo      = "Fruta"
$productoId = 100107
$producto  = "Otros"
$categoriaId = 100107002
$categoria = "Chirimoya"
$variedad  = "Cultivar IV Región"
$unidad    = "`$/bandeja 10 kilos"
$origen    = "Provincia de Limarí"
$kgUnidad  = 10

function Set-DataRow {
    param(
        [int]$Row,
        [double]$Fecha,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [double]$PrecioKg
    )

    $ws.Cells.Item($Row, 1).Value  = $mercadoId
    $ws.Cells.Item($Row, 2).Value  = $mercado
    $ws.Cells.Item($Row, 3).Value  = $region
    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 5).Value  = $codreg
    $ws.Cells.Item($Row, 6).Value  = $tipo
    $ws.Cells.Item($Row, 7).Value  = $productoId
    $ws.Cells.Item($Row, 8).Value  = $producto
    $ws.Cells.Item($Row, 9).Value  = $categoriaId
    $ws.Cells.Item($Row, 10).Value = $categoria
    $ws.Cells.Item($Row, 11).Value = $variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $unidad
    $ws.Cells.Item($Row, 18).Value = $origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $kgUnidad

    # Column D carries the date format already used throughout the sheet (style index 2).
    $ws.Cells.Item($Row, 4).NumberFormat = $ws.Cells.Item($Row - 1, 4).NumberFormat
}

# 1) Insert a new row above row 5 (pushes current rows 5-18 down to 6-19) and fill it in.
$ws.Rows.Item(5).Insert()
Set-DataRow 5 44421 "Especial" 30 35000 35000 35000 3500

# 2) Append two new rows at the bottom of the table (rows 20 and 21).
Set-DataRow 20 44461 "Especial" 150 30000 30000 30000 3000
Set-DataRow 21 44461 "Primera"  100 25000 25000 25000 2500
